$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the timestamp string in A2
$ws.Range("A2").Value = "2025-05-13 15:59:34"

# Update numeric metrics in row 2
$ws.Range("B2").Value = 16150
$ws.Range("C2").Value = 11659
$ws.Range("D2").Value = 72.19195046439629
$ws.Range("E2").Value = 2280
$ws.Range("F2").Value = 14.11764705882353
$ws.Range("G2").Value = 3096
$ws.Range("H2").Value = 19.1702786377709
$ws.Range("I2").Value = 9527
$ws.Range("J2").Value = 58.99071207430341
$ws.Range("K2").Value = 3044987.56
$ws.Range("L2").Value = 3527
$ws.Range("M2").Value = 21.83900928792569
$ws.Range("N2").Value = 1081584.63
$ws.Range("O2").Value = 4644
$ws.Range("P2").Value = 28.75541795665635
$ws.Range("Q2").Value = 481281.6
$ws.Range("R2").Value = 3518
$ws.Range("S2").Value = 21.78328173374613
$ws.Range("T2").Value = 3475
$ws.Range("U2").Value = 21.51702786377709
$ws.Range("V2").Value = 2428604.86
$ws.Range("W2").Value = 2302
$ws.Range("X2").Value = 14.25386996904025
$ws.Range("Y2").Value = 1408
$ws.Range("Z2").Value = 8.718266253869968
$ws.Range("AA2").Value = 135101.1
$ws.Range("AB2").Value = 803
$ws.Range("AC2").Value = 4.972136222910216
$ws.Range("AD2").Value = 456
$ws.Range("AE2").Value = 871
$ws.Range("AF2").Value = 1391
$ws.Range("AG2").Value = 16.77704194260486
$ws.Range("AH2").Value = 32.04562178072112
$ws.Range("AI2").Value = 51.17733627667403
$ws.Range("AJ2").Value = 1423477.08
$ws.Range("AK2").Value = 267238.87
$ws.Range("AL2").Value = 88992.39
$ws.Range("AM2").Value = 79.98372812030537
$ws.Range("AN2").Value = 15.01588007392268
$ws.Range("AO2").Value = 5.00039180577195
$ws.Range("AP2").Value = 46.56639839034205
$ws.Range("AQ2").Value = 207.1931755641167
$ws.Range("AR2").Value = 574.3230983949755
